$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_2_06_C")

# --- Title text updates (October -> November) ---
$ws.Range("A2").Value = "Useful Thermal Output, by Sector, 2006-November 2016 (Thousand Tons)"

# --- Insert a new row for "November" month data right after the October row (row 52) ---
$ws.Rows.Item(53).Insert()

# New row 53: November month row, formatted/styled like the other month rows (copy from row 52)
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 1455
$ws.Range("C53").Value = 39
$ws.Range("D53").Value = 1120
$ws.Range("E53").Value = 281
$ws.Range("F53").Value = 15

# --- Update "Annual Totals" rolling section (now shifted down by one row) ---
# Row 55: 2014
$ws.Range("B55").Value = 17095
$ws.Range("C55").Value = 408
$ws.Range("D55").Value = 13245
$ws.Range("E55").Value = 3253
$ws.Range("F55").Value = 189

# Row 56: 2015
$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 16971
$ws.Range("C56").Value = 415
$ws.Range("D56").Value = 13159
$ws.Range("E56").Value = 3203
$ws.Range("F56").Value = 193

# Row 57: 2016
$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 16463
$ws.Range("C57").Value = 425
$ws.Range("D57").Value = 12782
$ws.Range("E57").Value = 3096
$ws.Range("F57").Value = 159

# --- "Rolling 12 Months Ending in" header text update ---
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Row 59: 2015
$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 18537
$ws.Range("C59").Value = 451
$ws.Range("D59").Value = 14373
$ws.Range("E59").Value = 3502
$ws.Range("F59").Value = 211

# Row 60: 2016
$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 18108
$ws.Range("C60").Value = 462
$ws.Range("D60").Value = 14075
$ws.Range("E60").Value = 3395
$ws.Range("F60").Value = 177
